# Applies the edit described by the commit:
#  "Changed the way events work, removing the data class GameEvent."
# The visible textual content of the document does not change; Word's
# proofing engine (spell/grammar check) re-ran over the edited areas and
# left behind <w:proofErr/> markers, run splits were introduced at the
# proofing-mark boundaries, and the author pressed Enter in the middle of
# "The TimePowerupTest class - Tradition Breaker" (leaving the _GoBack
# bookmark behind) which split that paragraph into two list items. We
# reconstruct the final OOXML for each affected paragraph and push it in
# with Range.InsertXML, which is the only way to get proofErr/bookmark
# markup into the document via this object model.

$d = $word.ActiveDocument

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$apos = [char]0x2019

# ---------------------------------------------------------------------
# 1) Paragraph 8 - "The reason this class is a tradition break ..."
#    TimePowerupTest and the "@before" split get proofErr wrappers.
#    (Done first, bottom-up, so earlier paragraph indices below are not
#    disturbed by the paragraph-count change from edit #2.)
# ---------------------------------------------------------------------
$p8 = $d.Paragraphs(8).Range
$rng8 = $d.Range($p8.Start, $p8.End - 1)

$body8 = '<w:p>' + `
  '<w:r><w:t xml:space="preserve">The reason this class is a tradition break is that </w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">within the </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>TimePowerupTest</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r><w:t>o</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">nly </w:t></w:r>' + `
  '<w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">the </w:t></w:r>' + `
  '<w:r><w:t>@</w:t></w:r><w:proofErr w:type="gramEnd"/>' + `
  ('<w:r><w:t>before statement is implemented, the other test methods don' + $apos + 't really test anything.</w:t></w:r>') + `
  '</w:p>'

$rng8.InsertXML($pkgOpen + $body8 + $pkgClose)

# ---------------------------------------------------------------------
# 2) Paragraphs 6-7 - the "<< TODO" bookmark moves out of paragraph 6
#    and "The TimePowerupTest class - Tradition Breaker" is split into
#    two list paragraphs: "The " (carrying the _GoBack bookmark) and a
#    second one with the rest of the sentence (now with a proofErr
#    wrapper around TimePowerupTest).
# ---------------------------------------------------------------------
$p6 = $d.Paragraphs(6).Range
$p7 = $d.Paragraphs(7).Range
$rng67 = $d.Range($p6.Start, $p7.End - 1)

$pPr1 = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>'
$pPr0 = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>'

$body67 = ('<w:p>' + $pPr1 + '<w:r><w:t>&lt;&lt; TODO</w:t></w:r></w:p>') + `
  ('<w:p>' + $pPr0 + '<w:r><w:t xml:space="preserve">The </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>') + `
  ('<w:p>' + $pPr0 + '<w:r><w:t xml:space="preserve">The </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>TimePowerupTest</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> class - Tradition Breaker</w:t></w:r></w:p>')

$rng67.InsertXML($pkgOpen + $body67 + $pkgClose)

# ---------------------------------------------------------------------
# 3) Paragraph 5 - "The Resolution class is similar to the GameEvent
#    class ..." picks up proofErr wrappers around GameEvent, equals(,
#    and hashmap.
# ---------------------------------------------------------------------
$p5 = $d.Paragraphs(5).Range
$rng5 = $d.Range($p5.Start, $p5.End - 1)

$body5 = '<w:p>' + `
  '<w:r><w:t xml:space="preserve">The Resolution class is similar to the </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>GameEvent</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> cla</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">ss, it contains mainly a constructor, getters and setters. </w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">But it has additional </w:t></w:r>' + `
  '<w:proofErr w:type="gramStart"/><w:r><w:t>equals(</w:t></w:r><w:proofErr w:type="gramEnd"/>' + `
  '<w:r><w:t xml:space="preserve">) and </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>hashmap</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t>() methods, which makes the class look like a normal class. But these methods are just methods overridden from the Object superclass. Therefore this class is indeed a data class.</w:t></w:r>' + `
  '</w:p>'

$rng5.InsertXML($pkgOpen + $body5 + $pkgClose)
